$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 31; $r -le 91; $r++) {
    $ws.Cells.Item($r, 8).Value = 0    # Column H
    $ws.Cells.Item($r, 10).Value = 1   # Column J
}

for ($r = 92; $r -le 138; $r++) {
    $ws.Cells.Item($r, 8).Value = 0      # Column H
    $ws.Cells.Item($r, 10).Value = 0.2   # Column J
}

for ($r = 139; $r -le 176; $r++) {
    $ws.Cells.Item($r, 8).Value = 0    # Column H only
}
